$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.255.53"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").Value = "1.863.98"
$ws.Range("E3").Value = "  -1.15%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "0.7054"
$ws.Range("E5").Value = "  -1.17%  "
$ws.Range("D6").Value = "242.51"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "0.07827"
$ws.Range("E8").Value = "  -3.18%  "
$ws.Range("D9").Value = "0.3110"
$ws.Range("E9").Value = "  -0.68%  "
$ws.Range("D10").Value = "24.27"
$ws.Range("E10").Value = "  -4.02%  "
$ws.Range("D11").Value = "0.07992"
$ws.Range("E11").Value = "  -4.36%  "
$ws.Range("D12").Value = "1.893.35"
$ws.Range("E12").Value = "  +0.79%  "
$ws.Range("D13").Value = "5.180"
$ws.Range("E13").Value = "  -1.27%  "
$ws.Range("D14").Value = "93.57"
$ws.Range("E14").Value = "  +1.37%  "
$ws.Range("D15").Value = "0.6946"
$ws.Range("E15").Value = "  -3.72%  "
$ws.Range("D16").Value = "6.358"
$ws.Range("E16").Value = "  +0.95%  "
$ws.Range("D17").Value = "29.413.47"
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D18").Value = "0.000008278"
$ws.Range("E18").Value = "  -2.24%  "
$ws.Range("D19").Value = "251.55"
$ws.Range("E19").Value = "  +4.22%  "
$ws.Range("D20").Value = "2.175.72"
$ws.Range("E20").Value = "  +2.11%  "
$ws.Range("D21").Value = "13.11"
$ws.Range("E21").Value = "  -1.08%  "
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").Value = "7.524"
$ws.Range("E23").Value = "  -3.82%  "
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("E25").Value = "  -2.34%  "
$ws.Range("D26").Value = "8.994"
$ws.Range("E26").Value = "  -0.96%  "
$ws.Range("D27").Value = "159.67"
$ws.Range("E27").Value = "  -2.28%  "
$ws.Range("D28").Value = "18.68"
$ws.Range("E28").Value = "  +0.53%  "
$ws.Range("D29").Value = "1.499"
$ws.Range("E29").Value = "  -0.50%  "
$ws.Range("D30").Value = "4.273"
$ws.Range("E30").Value = "  -1.54%  "
$ws.Range("D31").Value = "4.269"
$ws.Range("E31").Value = "  -3.65%  "
$ws.Range("D32").Value = "1.214"
$ws.Range("E32").Value = "  +0.42%  "
$ws.Range("D33").Value = "0.05257"
$ws.Range("D34").Value = "1.892"
$ws.Range("E34").Value = "  -3.20%  "
$ws.Range("D35").Value = "0.7455"
$ws.Range("E35").Value = "  -0.69%  "
$ws.Range("D36").Value = "1.155"
$ws.Range("E36").Value = "  -2.26%  "
$ws.Range("D37").Value = "2.704"
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("D38").Value = "0.01858"
$ws.Range("E38").Value = "  -1.23%  "
$ws.Range("D39").Value = "1.241.94"
$ws.Range("E39").Value = "  -3.24%  "
$ws.Range("D40").Value = "2.746"
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").Value = "6.207"
$ws.Range("E41").Value = "  -5.50%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "0.9011"
$ws.Range("E42").Value = "  +0.95%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "111.07"
$ws.Range("E43").Value = "  +0.57%  "
$ws.Range("D44").Value = "71.82"
$ws.Range("E44").Value = "  -2.33%  "
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").Value = "1.000"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "2.080.99"
$ws.Range("E46").Value = "  +2.61%  "
$ws.Range("D47").Value = "0.00000000123"
$ws.Range("E47").Value = "  -4.40%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "0.5201"
$ws.Range("E48").Value = "  -0.23%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "1.790"
$ws.Range("E49").Value = "  -0.94%  "
$ws.Range("D50").Value = "9.384"
$ws.Range("E50").Value = "  -1.20%  "
$ws.Range("D51").Value = "1.014"
$ws.Range("E51").Value = "  +0.96%  "
